$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UploadedFileSource")

# Insert a new column before column E (shifts input/output/type/output_name/database_id/id right by one)
$ws.Columns.Item(5).Insert()

# Set header for the newly inserted column E
$ws.Cells.Item(1, 5).Value = "coordinates"
